# Update the "Förändrad" (Changed) date column C for all data rows (2-38)
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C38").Value = 45174
